$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move rows down to create gaps: insert a blank row before row 3 (shifts old row3 data to row4)
$ws.Rows("3").Insert()
# Now insert another blank row before what is now row 5 (old row4 data), shifting it to row 6
$ws.Rows("5").Insert()
# Now insert another blank row before what is now row 7 (old row5 data), shifting it to row 8
$ws.Rows("7").Insert()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:F8"), 0, 1)
$tbl.Name = "Table4"
$tbl.TableStyle = "TableStyleMedium13"

$ws.Range("I9:I10").Select()
$excel.ActiveCell = $ws.Range("I10")

Write-Host $tbl.Name
